# Export with no is_pref and no lev distance
#
# The speaker/variant list is re-exported: entries are reordered, two
# ids are corrected to match their variant text, a new "#2"/"2" row is
# added, the is_prefered ("x") column is cleared for every row, and a
# "#prudent"/"Prudent" row is appended at the end (row 24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://www.dbnl.org/tekst/baud004edip01_01"

# New ordering of (id, speaker_variant) for rows 2..24
$data = @(
    @("#edipes", "Edipes"),
    @("#vridies", "Vridies"),
    @("#iocaste", "Iocaste"),
    @("#hemon", "Hemon"),
    @("#antigone", "Antigone"),
    @("#polenies", "Polenies"),
    @("#macolph", "Macolph"),
    @("#creon", "Creon"),
    @("#bode", "Bode"),
    @("#broersh", "Broersh"),
    @("#eteocle", "Eteocle"),
    @("#ismene", "Ismene"),
    @("#1.-soldaet", "1. Soldaet"),
    @("#koddig", "Koddig"),
    @("#antigon", "Antigon"),
    @("#cerbo", "Cerbo"),
    @("#echo", "Echo"),
    @("#sitolus", "Sitolus"),
    @("#galo", "Galo"),
    @("#macolp", "Macolp"),
    @("#2", "2"),
    @("#antetus", "Antetus"),
    @("#prudent", "Prudent")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $url
    $ws.Cells.Item($row, 2).Value = $entry[0]

    $variant = $entry[1]
    $cVariant = $ws.Cells.Item($row, 3)
    if ($variant -eq "2") {
        # Keep this purely-numeric-looking variant stored as text, not a number.
        $cVariant.NumberFormat = "@"
    }
    $cVariant.Value = $variant

    $ws.Cells.Item($row, 4).ClearContents()
    $row = $row + 1
}

# Row 24 is brand new, so its D:H cells (present-but-empty in the original
# sheet, like every other data row) need to be created explicitly: touch a
# lightweight formatting property to force the cell node to exist, then
# clear any content so the cell stays empty.
$ws.Range("D24:H24").WrapText = $false
$ws.Range("D24:H24").ClearContents()
